# [Closed] escape first n rows, sheet name or number #1
# Add a second worksheet ("Sheet2") that mirrors Sheet1's "option" data,
# with a title/subtitle header block above the normal header row.

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

# Create Sheet2 as a full copy of Sheet1 (keeps column widths, page setup,
# phonetic settings, etc. identical to Sheet1) and place it right after Sheet1.
$sheet1.Copy($null, $sheet1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Sheet2"

# --- Sheet1: only the selected cell/view changes ---
$sheet1.Range("C7").Select()

# --- Sheet2: insert two new rows on top for a title + subtitle block ---
$ws2.Range("A1:A2").EntireRow.Insert()

$ws2.Range("A1").Value = "title"
$ws2.Range("A2").Value = "subtitle1"
$ws2.Range("C2").Value = "subtitle2"

# Merge & center the title / subtitle cells
$ws2.Range("A1:C1").Merge()
$ws2.Range("A2:B2").Merge()
$ws2.Range("A1:C2").HorizontalAlignment = -4108
$ws2.Range("A1:C2").VerticalAlignment = -4108

# Replace the option data (columns B & C, rows 4-11) with the "_option_" set
$ws2.Range("B4").Value = "1_option_name_xlsx"
$ws2.Range("C4").Value = "1_option_xlsx@email.com"
$ws2.Range("B5").Value = "2_option_name_xlsx"
$ws2.Range("C5").Value = "2_option_xlsx@email.com"
$ws2.Range("B6").Value = "3_option_name_xlsx"
$ws2.Range("C6").Value = "3_option_xlsx@email.com"
$ws2.Range("C7").Value = "4_option_xlsx@email.com"
$ws2.Range("B8").Value = "5_option_name_xlsx"
$ws2.Range("C8").Value = "5_option_xlsx@email.com"
$ws2.Range("B9").Value = "6_option_name_xlsx"
$ws2.Range("B10").Value = "7_option_name_xlsx"
$ws2.Range("C10").Value = "7_option_xlsx@email.com"
$ws2.Range("B11").Value = "8_option_name_xlsx"
$ws2.Range("C11").Value = "8_option_xlsx@email.com"

# Sheet2 ends up the active sheet/tab with C9 selected
$ws2.Activate()
$ws2.Range("C9").Select()
